$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix 1: correct a stray comma to a period in a proveedor name ---
$oldName = "IZAGUIRRE CARLOS MARIA, MOREND MARIA ELENA Y MOREND MARIA TERESA"
$newName = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Cells.Replace($oldName, $newName)

# --- Fix 2: re-format "Importe" amounts from "1.234,56" (ES-style) text
#            to "1234.56" (plain decimal-point) text, column H, rows 2-80.
#            A helper sheet computes each literal via TEXT(1,"...") -- a
#            formula that always yields a Text/string result -- then a
#            Copy + PasteSpecial(xlPasteValues) round-trip drops the
#            formula and keeps only the resulting string, landing in the
#            destination cells as plain shared-string text (never
#            reinterpreted as a number) without touching their style.
$values = @(
    "1730.00",
    "37498.60",
    "35000.00",
    "3225.84",
    "24220.00",
    "96496.84",
    "8296.75",
    "1240.50",
    "766.00",
    "5158.16",
    "231.90",
    "5736.22",
    "89.60",
    "3585.00",
    "1598.04",
    "1140.00",
    "360.00",
    "18.12",
    "4352.00",
    "1250.00",
    "544.00",
    "14352.00",
    "1472.00",
    "3168.00",
    "7904.00",
    "9000.00",
    "13036.50",
    "392.04",
    "426.00",
    "1280.00",
    "903.50",
    "118.60",
    "451425.96",
    "138.50",
    "198.40",
    "450.00",
    "1770.00",
    "780.90",
    "2492.00",
    "4443.99",
    "3135.00",
    "4400.00",
    "1570.00",
    "250.00",
    "500.00",
    "890.00",
    "1013.26",
    "12491.15",
    "500.00",
    "1000.00",
    "1400.00",
    "150.00",
    "4770.31",
    "570.00",
    "386.75",
    "3630.00",
    "250.00",
    "628.00",
    "8563.98",
    "600.00",
    "1350.00",
    "1326.00",
    "300.00",
    "2000.00",
    "4663.00",
    "662.58",
    "170.00",
    "180.00",
    "168.22",
    "2285.00",
    "1118.00",
    "225.00",
    "475.00",
    "805.50",
    "966.00",
    "3179.88",
    "7007.76",
    "350716.53",
    "991.26"
)

$helper = $wb.Worksheets.Add($null, $ws)
for ($i = 0; $i -lt $values.Length; $i++) {
    $literal = $values[$i].Replace('"', '""')
    $helper.Cells.Item($i + 1, 1).Formula = '=TEXT(1,"""' + $literal + '""")'
}

$firstRow = 2
$lastRow = $firstRow + $values.Length - 1
$srcRange = $helper.Range($helper.Cells.Item(1, 1), $helper.Cells.Item($values.Length, 1))
$dstRange = $ws.Range("H" + $firstRow + ":H" + $lastRow)

$srcRange.Copy()
$dstRange.PasteSpecial(-4163)

$helper.Delete()
